# Apply "IBG Cao Scripts" edit:
#  - Rename Sheet1 -> CFR, Sheet2 -> CFR_Ibg
#  - Duplicate the CFR data (header + data row) onto CFR_Ibg
#  - Copy CFR's column widths onto CFR_Ibg
#  - Make CFR_Ibg the active/selected tab, with B19 selected on CFR
#    and D16 selected on CFR_Ibg

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Rename the sheets
$ws1.Name = "CFR"
$ws2.Name = "CFR_Ibg"

# Duplicate the data from CFR onto CFR_Ibg, preserving formats/number styles
$ws1.Range("A1:L2").Copy($ws2.Range("A1"))

# Match column widths to CFR's
for ($i = 1; $i -le 12; $i++) {
    $ws2.Columns.Item($i).ColumnWidth = $ws1.Columns.Item($i).ColumnWidth
}

# Update the selection on CFR (no longer the active tab)
$ws1.Range("B19").Select() | Out-Null

# Make CFR_Ibg the active tab with D16 selected
$ws2.Activate() | Out-Null
$ws2.Range("D16").Select() | Out-Null
